$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 142868.86
$ws.Range("I11").Value = 142868.86
$ws.Range("K11").Value = 142868.86
$ws.Range("M11").Value = -142728.86
$ws.Range("H15").Value = 1084.9
$ws.Range("I15").Value = 1084.9
$ws.Range("K15").Value = 3254.7
$ws.Range("M15").Value = -3085.7
$ws.Range("H41").Value = 812.63635
$ws.Range("I41").Value = 1293.5385
$ws.Range("J41").Value = 500.05
$ws.Range("K41").Value = 1293.5385
$ws.Range("L41").Value = 500.05
$ws.Range("M41").Value = -853.5385000000001
$ws.Range("N41").Value = -1380.05
$ws.Range("H62").Value = 2328.5715
$ws.Range("I62").Value = 2328.5715
$ws.Range("K62").Value = 2328.5715
$ws.Range("M62").Value = -1704.5715
$ws.Range("H65").Value = 2328.5715
$ws.Range("I65").Value = 2328.5715
$ws.Range("K65").Value = 11642.8575
$ws.Range("M65").Value = -8522.8575
$ws.Range("H106").Value = 2993
$ws.Range("I106").Value = 4233.1665
$ws.Range("K106").Value = 4233.1665
$ws.Range("M106").Value = -3602.1665
$ws.Range("H113").Value = 85012.836
$ws.Range("I113").Value = 251048.75
$ws.Range("J113").Value = 1994.875
$ws.Range("K113").Value = 251048.75
$ws.Range("L113").Value = 1994.875
$ws.Range("M113").Value = -247794.75
$ws.Range("N113").Value = -8502.875
$ws.Range("H129").Value = 2392.972
$ws.Range("J129").Value = 1142.585
$ws.Range("L129").Value = 3427.755
$ws.Range("N129").Value = -13427.755
$ws.Range("H138").Value = 1646.9512
$ws.Range("I138").Value = 1243.8788
$ws.Range("J138").Value = 3309.625
$ws.Range("K138").Value = 3731.6364
$ws.Range("L138").Value = 9928.875
$ws.Range("M138").Value = 1408.3636
$ws.Range("N138").Value = -20208.875

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4755.2
$ws.Range("I32").Value = 3658.3635
$ws.Range("J32").Value = 12798.667
$ws.Range("K32").Value = 3658.3635
$ws.Range("L32").Value = 12798.667
$ws.Range("M32").Value = -3371.3635
$ws.Range("N32").Value = -13372.667
$ws.Range("H74").Value = 443.36365
$ws.Range("I74").Value = 324.45456
$ws.Range("K74").Value = 324.45456
$ws.Range("M74").Value = 549.54544
$ws.Range("H77").Value = 443.36365
$ws.Range("I77").Value = 324.45456
$ws.Range("K77").Value = 1622.2728
$ws.Range("M77").Value = 2745.7272
$ws.Range("H122").Value = 1427.3
$ws.Range("I122").Value = 1412.1111
$ws.Range("K122").Value = 4236.3333
$ws.Range("M122").Value = -1786.3333

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 13000
$ws.Range("J44").Value = 13000
$ws.Range("L44").Value = 13000
$ws.Range("N44").Value = -13884
$ws.Range("H132").Value = 2712.8147
$ws.Range("I132").Value = 2568.1462
$ws.Range("J132").Value = 3169.077
$ws.Range("K132").Value = 7704.4386
$ws.Range("L132").Value = 9507.231
$ws.Range("M132").Value = -5174.4386
$ws.Range("N132").Value = -14567.231

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 921.5952
$ws.Range("I5").Value = 500.34784
$ws.Range("J5").Value = 1431.5264
$ws.Range("K5").Value = 1501.04352
$ws.Range("L5").Value = 4294.5792
$ws.Range("M5").Value = -1389.04352
$ws.Range("N5").Value = -4518.5792
$ws.Range("H23").Value = 583.4706
$ws.Range("I23").Value = 290
$ws.Range("J23").Value = 673.7692
$ws.Range("K23").Value = 870
$ws.Range("L23").Value = 2021.3076
$ws.Range("M23").Value = -635
$ws.Range("N23").Value = -2491.3076
$ws.Range("H37").Value = 647474.4
$ws.Range("J37").Value = 647474.4
$ws.Range("L37").Value = 1942423.2
$ws.Range("N37").Value = -1942647.2
$ws.Range("H41").Value = 793.3333
$ws.Range("I41").Value = 793.3333
$ws.Range("K41").Value = 2379.9999
$ws.Range("M41").Value = -2041.9999
$ws.Range("H51").Value = 3300
$ws.Range("I51").Value = 1066.6666
$ws.Range("J51").Value = 4257.143
$ws.Range("K51").Value = 3199.9998
$ws.Range("L51").Value = 12771.429
$ws.Range("M51").Value = -2739.9998
$ws.Range("N51").Value = -13691.429
$ws.Range("H113").Value = 555.4516
$ws.Range("J113").Value = 555.8095
$ws.Range("L113").Value = 1667.4285
$ws.Range("N113").Value = -6007.4285
$ws.Range("H131").Value = 1333.4615
$ws.Range("I131").Value = 914.8570999999999
$ws.Range("J131").Value = 1368.3452
$ws.Range("K131").Value = 2744.5713
$ws.Range("L131").Value = 4105.0356
$ws.Range("M131").Value = 2295.4287
$ws.Range("N131").Value = -14185.0356
$ws.Range("H135").Value = 921.5952
$ws.Range("I135").Value = 500.34784
$ws.Range("J135").Value = 1431.5264
$ws.Range("K135").Value = 4503.130560000001
$ws.Range("L135").Value = 12883.7376
$ws.Range("M135").Value = -1968.130560000001
$ws.Range("N135").Value = -17953.7376

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("H135").Value = 27631.295
$ws.Range("J135").Value = 27631.295
$ws.Range("L135").Value = 27631.295
$ws.Range("N135").Value = -37771.295
$ws.Range("H139").Value = 44645
$ws.Range("J139").Value = 44645
$ws.Range("L139").Value = 44645
$ws.Range("N139").Value = -54925
$ws.Range("N104").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 69065.336
$ws.Range("I40").Value = 252120
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 252120
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -251984
$ws.Range("N40").Value = -2772
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("H132").Value = 3061.394
$ws.Range("I132").Value = 3211.1035
$ws.Range("J132").Value = 1976
$ws.Range("K132").Value = 9633.3105
$ws.Range("L132").Value = 5928
$ws.Range("M132").Value = -7103.3105
$ws.Range("N132").Value = -10988
$ws.Range("H135").Value = 35675
$ws.Range("J135").Value = 35675
$ws.Range("L135").Value = 35675
$ws.Range("N135").Value = -45815
$ws.Range("H136").Value = 1135.6086
$ws.Range("I136").Value = 1085.5897
$ws.Range("J136").Value = 1414.2858
$ws.Range("K136").Value = 3256.7691
$ws.Range("L136").Value = 4242.857400000001
$ws.Range("M136").Value = -706.7691
$ws.Range("N136").Value = -9342.857400000001
$ws.Range("M75").ClearContents()
$ws.Range("M78").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 35000
$ws.Range("J127").Value = 35000
$ws.Range("L127").Value = 35000
$ws.Range("N127").Value = -44920
$ws.Range("H133").Value = 39500
$ws.Range("J133").Value = 39500
$ws.Range("L133").Value = 39500
$ws.Range("N133").Value = -49620
$ws.Range("H135").Value = 48999.75
$ws.Range("J135").Value = 48999.75
$ws.Range("L135").Value = 48999.75
$ws.Range("N135").Value = -59139.75
